# Base bar plot for Fig 1 panel B - add oro_type factor level aesthetics
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 5-11: oro_type factor aesthetics ---
# Column A (rows 5-11): variable name "oro_type"
for ($r = 5; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = "oro_type"
}

# Column B (rows 5-11): labels
$bVals = @("Marine renewable energy", "CO2 removal or storage", "Increase efficiency", "Human assisted evolution", "Conservation", "Built infrastructure & technology", "Socio-institutional")
for ($i = 0; $i -lt 7; $i++) {
    $ws.Cells.Item(5 + $i, 2).Value = $bVals[$i]
}

# Column C (rows 5-11): levels (note row 8 uses hyphenated spelling)
$cVals = @("Marine renewable energy", "CO2 removal or storage", "Increase efficiency", "Human-assisted evolution", "Conservation", "Built infrastructure & technology", "Socio-institutional")
for ($i = 0; $i -lt 7; $i++) {
    $ws.Cells.Item(5 + $i, 3).Value = $cVals[$i]
}

# Column D (rows 5-11): order numbers 1-7
for ($i = 0; $i -lt 7; $i++) {
    $ws.Cells.Item(5 + $i, 4).Value = $i + 1
}

# Column E (rows 5-11): colour formulas referencing the branch-level colours
$eRefs = @("E2", "E2", "E2", "E3", "E3", "E4", "E4")
for ($i = 0; $i -lt 7; $i++) {
    $ws.Cells.Item(5 + $i, 5).Formula = "=" + $eRefs[$i]
}

# --- Formatting / layout tweaks ---
# Widen column C to fit the new longer labels
$ws.Columns.Item(3).ColumnWidth = 15.619791666666666

# Move the active selection as recorded in the saved file
$ws.Range("E12").Select() | Out-Null

# Page setup: portrait orientation
$ws.PageSetup.Orientation = 1

Write-Host "done"
